$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "SUBSCRIPCIÓN" column header in G1, reusing the same (shaded)
# header style that the existing header cells (e.g. F1) already have.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "SUBSCRIPCIÓN"

# Column F (previously the last data column, width 17.5 with outlineLevel 1
# like the rest of the grouped columns) keeps the width/outline that column
# E already has, and the new column G takes the width that column F used to
# have (17, not part of the outline group).
$ws.Columns("F").ColumnWidth = 16.666666666666668
$ws.Columns("F").OutlineLevel = $ws.Columns("E").OutlineLevel
$ws.Columns("G").ColumnWidth = 16.166666666666668

# The active selection moves to G2, ready for the first data row under the
# new header.
$unused = $ws.Range("G2").Select()

# Extend the autofilter so it covers the new column too.
$ws.AutoFilterMode = $false
$unused = $ws.Range("A1:G1").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$n = $wb.Names.Item("Data!_FilterDatabase")
$n.RefersTo = "=Data!`$A`$1:`$G`$1"
